# Weekly update: insert a new price-report row for "Haba" at row 56,
# pushing the existing rows 56-114 down to 57-115 (new dimension A1:R115).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 56 (Excel shifts rows 56.. down by one
# and copies formatting from the row above, matching the date style used
# in column D).
$ws.Rows.Item(56).Insert()

# Populate the newly inserted row 56 with the new record.
$ws.Range("A56").Value = 3
$ws.Range("B56").Value = "Femacal de La Calera"
$ws.Range("C56").Value = "Coquimbo"
$ws.Range("D56").Value = 44546
$ws.Range("E56").Value = 5
$ws.Range("F56").Value = 100112026
$ws.Range("G56").Value = "Haba"
$ws.Range("H56").Value = "Sin especificar"
$ws.Range("I56").Value = "Primera"
$ws.Range("J56").Value = 95
$ws.Range("K56").Value = 8000
$ws.Range("L56").Value = 8500
$ws.Range("M56").Value = 8237
$ws.Range("N56").Value = "`$/saco 25 kilos"
$ws.Range("O56").Value = "Provincia de Quillota"
$ws.Range("P56").Value = 329
$ws.Range("Q56").Value = 25
$ws.Range("R56").Value = "Hortaliza"
